$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column BF holds the game "Date" as text (e.g. "4-20-2012-13"). The NBA
# stats export mis-formatted it; fix it to ISO form "2013-04-20" for every
# data row (BF2:BF31), keeping the value as plain text (not an Excel date
# serial) and leaving the cell's style untouched.
for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Cells.Item($r, 58)
    $cell.NumberFormat = "@"        # force text interpretation while we set it
    $cell.Value = "2013-04-20"
    $cell.Style = "Normal"          # restore the original (default) cell style
}
